$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Footer 1 (default footer) - Pearson logo: image1.png -> image2.png
$ftr1 = $sec.Footers(1)
if ($ftr1.Exists -and $ftr1.Range.InlineShapes.Count -ge 1) {
    $ftr1.Range.InlineShapes(1).Name = "image2.png"
}
Write-Output "Updated footer1 Pearson logo name"

# Footer 2 (first-page footer) - Pearson logo: image1.png -> image2.png
$ftr2 = $sec.Footers(2)
if ($ftr2.Exists -and $ftr2.Range.InlineShapes.Count -ge 1) {
    $ftr2.Range.InlineShapes(1).Name = "image2.png"
}
Write-Output "Updated footer2 Pearson logo name"

# Header 2 (first-page header) - BTec logo: image2.jpg -> image1.jpg
$hdr2 = $sec.Headers(2)
if ($hdr2.Exists -and $hdr2.Range.InlineShapes.Count -ge 1) {
    $hdr2.Range.InlineShapes(1).Name = "image1.jpg"
}
Write-Output "Updated header2 BTec logo name"
